$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H125").Value = 835.0526
$ws.Range("I125").Value = 458.8889
$ws.Range("J125").Value = 1173.6
$ws.Range("K125").Value = 4130.0001
$ws.Range("L125").Value = 10562.4
$ws.Range("M125").Value = -1670.0001
$ws.Range("N125").Value = -15482.4
$ws.Range("H132").Value = 3117.8948
$ws.Range("I132").Value = 2773.3655
$ws.Range("J132").Value = 6701
$ws.Range("K132").Value = 8320.0965
$ws.Range("L132").Value = 20103
$ws.Range("M132").Value = -5790.0965
$ws.Range("N132").Value = -25163
$ws.Range("H135").Value = 1721.0834
$ws.Range("I135").Value = 1783.3636
$ws.Range("J135").Value = 1036
$ws.Range("K135").Value = 16050.2724
$ws.Range("L135").Value = 9324
$ws.Range("M135").Value = -13515.2724
$ws.Range("N135").Value = -14394
$ws.Range("H137").Value = 1625.4445
$ws.Range("J137").Value = 1825.5714
$ws.Range("L137").Value = 5476.7142
$ws.Range("N137").Value = -10576.7142
$ws.Range("H138").Value = 4651.564
$ws.Range("J138").Value = 4755.3228
$ws.Range("L138").Value = 14265.9684
$ws.Range("N138").Value = -24545.9684
$ws.Range("H141").Value = 4266.125
$ws.Range("I141").Value = 2012.6428
$ws.Range("J141").Value = 7421
$ws.Range("K141").Value = 6037.928400000001
$ws.Range("L141").Value = 22263
$ws.Range("M141").Value = -857.9284000000007
$ws.Range("N141").Value = -32623

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2079.0715
$ws.Range("I45").Value = 1545.4546
$ws.Range("K45").Value = 1545.4546
$ws.Range("M45").Value = -1168.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("H94").Value = 722.24
$ws.Range("I94").Value = 679.6923
$ws.Range("J94").Value = 768.3333
$ws.Range("K94").Value = 679.6923
$ws.Range("L94").Value = 768.3333
$ws.Range("M94").Value = -228.6923
$ws.Range("N94").Value = -1670.3333
$ws.Range("H105").Value = 2619.9
$ws.Range("I105").Value = 2281.8
$ws.Range("J105").Value = 2958
$ws.Range("K105").Value = 2281.8
$ws.Range("L105").Value = 2958
$ws.Range("M105").Value = -534.8000000000002
$ws.Range("N105").Value = -6452

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 30000
$ws.Range("I4").Value = 30000
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 30000
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -29888
$ws.Range("N4").ClearContents()
$ws.Range("H7").Value = 107.166664
$ws.Range("I7").Value = 76.22221999999999
$ws.Range("J7").Value = 200
$ws.Range("K7").Value = 76.22221999999999
$ws.Range("L7").Value = 200
$ws.Range("M7").Value = 36.77778000000001
$ws.Range("N7").Value = -426
$ws.Range("H31").Value = 1209.9183
$ws.Range("I31").Value = 1365.5652
$ws.Range("K31").Value = 1365.5652
$ws.Range("M31").Value = -1070.5652
$ws.Range("H34").Value = 1209.9183
$ws.Range("I34").Value = 1365.5652
$ws.Range("K34").Value = 1365.5652
$ws.Range("M34").Value = -1163.5652
$ws.Range("H58").Value = 1189.871
$ws.Range("I58").Value = 997.53656
$ws.Range("J58").Value = 1565.381
$ws.Range("K58").Value = 997.53656
$ws.Range("L58").Value = 1565.381
$ws.Range("M58").Value = -794.53656
$ws.Range("N58").Value = -1971.381
$ws.Range("H99").Value = 2000
$ws.Range("I99").Value = 2000
$ws.Range("J99").Value = 2000
$ws.Range("K99").Value = 2000
$ws.Range("L99").Value = 2000
$ws.Range("M99").Value = -502
$ws.Range("N99").Value = -4996
$ws.Range("H126").Value = 2000
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 2000
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 6000
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -10940
$ws.Range("H132").Value = 11113752
$ws.Range("I132").Value = 2282.818
$ws.Range("J132").Value = 41670292
$ws.Range("K132").Value = 6848.454000000001
$ws.Range("L132").Value = 125010876
$ws.Range("M132").Value = -4318.454000000001
$ws.Range("N132").Value = -125015936
$ws.Range("H136").Value = 1189.871
$ws.Range("I136").Value = 997.53656
$ws.Range("J136").Value = 1565.381
$ws.Range("K136").Value = 2992.60968
$ws.Range("L136").Value = 4696.143
$ws.Range("M136").Value = -442.60968
$ws.Range("N136").Value = -9796.143

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 12728571
$ws.Range("I4").Value = 4001428.5
$ws.Range("J4").Value = 100000000
$ws.Range("K4").Value = 12004285.5
$ws.Range("L4").Value = 300000000
$ws.Range("M4").Value = -12004173.5
$ws.Range("N4").Value = -300000224
$ws.Range("H45").Value = 924.6
$ws.Range("I45").Value = 200
$ws.Range("J45").Value = 1105.75
$ws.Range("K45").Value = 600
$ws.Range("L45").Value = 3317.25
$ws.Range("M45").Value = -68
$ws.Range("N45").Value = -4381.25
$ws.Range("H81").Value = 6812.5
$ws.Range("I81").Value = 2017
$ws.Range("K81").Value = 6051
$ws.Range("M81").Value = -4928
$ws.Range("H84").Value = 6812.5
$ws.Range("I84").Value = 2017
$ws.Range("K84").Value = 18153
$ws.Range("M84").Value = -12537
$ws.Range("H115").Value = 6004.8335
$ws.Range("I115").Value = 6019
$ws.Range("J115").Value = 5997.75
$ws.Range("K115").Value = 18057
$ws.Range("L115").Value = 17993.25
$ws.Range("M115").Value = -16882
$ws.Range("N115").Value = -20343.25
$ws.Range("H131").Value = 976.0526
$ws.Range("I131").Value = 657
$ws.Range("J131").Value = 1090
$ws.Range("K131").Value = 1971
$ws.Range("L131").Value = 3270
$ws.Range("M131").Value = 3069
$ws.Range("N131").Value = -13350

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 4333.3335
$ws.Range("I5").Value = 3000
$ws.Range("K5").Value = 3000
$ws.Range("M5").Value = -2888
$ws.Range("H11").Value = 5035002
$ws.Range("I11").Value = 10000000
$ws.Range("K11").Value = 10000000
$ws.Range("M11").Value = -9999861
$ws.Range("H12").Value = 52753.25
$ws.Range("I12").Value = 1001
$ws.Range("K12").Value = 1001
$ws.Range("M12").Value = -861
$ws.Range("H93").Value = 29582.5
$ws.Range("J93").Value = 29582.5
$ws.Range("L93").Value = 29582.5
$ws.Range("N93").Value = -33326.5
$ws.Range("H107").Value = 837.6
$ws.Range("I107").Value = 633.7143
$ws.Range("J107").Value = 1016
$ws.Range("K107").Value = 633.7143
$ws.Range("L107").Value = 1016
$ws.Range("M107").Value = 1286.2857
$ws.Range("N107").Value = -4856
$ws.Range("H113").Value = 1789.1666
$ws.Range("I113").Value = 1430.5
$ws.Range("J113").Value = 2506.5
$ws.Range("K113").Value = 1430.5
$ws.Range("L113").Value = 2506.5
$ws.Range("M113").Value = 739.5
$ws.Range("N113").Value = -6846.5
$ws.Range("H132").Value = 4959.8
$ws.Range("I132").Value = 5266.6665
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 15799.9995
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -13269.9995
$ws.Range("N132").Value = -18558.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H9").Value = 1255.1666
$ws.Range("I9").Value = 473.33334
$ws.Range("K9").Value = 473.33334
$ws.Range("M9").Value = -249.33334
$ws.Range("H61").Value = 3674.8125
$ws.Range("I61").Value = 3621.2
$ws.Range("J61").Value = 3764.1667
$ws.Range("K61").Value = 3621.2
$ws.Range("L61").Value = 3764.1667
$ws.Range("M61").Value = -3419.2
$ws.Range("N61").Value = -4168.1667
$ws.Range("H82").Value = 1777.2142
$ws.Range("I82").Value = 1156.4
$ws.Range("J82").Value = 2122.111
$ws.Range("K82").Value = 1156.4
$ws.Range("L82").Value = 2122.111
$ws.Range("M82").Value = -795.4000000000001
$ws.Range("N82").Value = -2844.111
$ws.Range("H85").Value = 1777.2142
$ws.Range("I85").Value = 1156.4
$ws.Range("J85").Value = 2122.111
$ws.Range("K85").Value = 1156.4
$ws.Range("L85").Value = 2122.111
$ws.Range("M85").Value = 91.59999999999991
$ws.Range("N85").Value = -4618.111
$ws.Range("H113").Value = 3674.8125
$ws.Range("I113").Value = 3621.2
$ws.Range("J113").Value = 3764.1667
$ws.Range("K113").Value = 3621.2
$ws.Range("L113").Value = 3764.1667
$ws.Range("M113").Value = -1451.2
$ws.Range("N113").Value = -8104.1667
$ws.Range("H123").Value = 65214.5
$ws.Range("J123").Value = 65214.5
$ws.Range("L123").Value = 65214.5
$ws.Range("N123").Value = -75014.5
$ws.Range("H132").Value = 3641.5588
$ws.Range("I132").Value = 3083.889
$ws.Range("J132").Value = 4268.9375
$ws.Range("K132").Value = 9251.667000000001
$ws.Range("L132").Value = 12806.8125
$ws.Range("M132").Value = -6721.667000000001
$ws.Range("N132").Value = -17866.8125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 255250000
$ws.Range("J2").Value = 500500000
$ws.Range("L2").Value = 500500000
$ws.Range("N2").Value = -500500224
$ws.Range("H17").Value = 37499.75
$ws.Range("I17").Value = 4995
$ws.Range("J17").Value = 48334.668
$ws.Range("K17").Value = 4995
$ws.Range("L17").Value = 48334.668
$ws.Range("M17").Value = -4823
$ws.Range("N17").Value = -48678.668
$ws.Range("H64").Value = 99999.5
$ws.Range("J64").Value = 99999.5
$ws.Range("L64").Value = 99999.5
$ws.Range("N64").Value = -100495.5
$ws.Range("H67").Value = 99999.5
$ws.Range("J67").Value = 99999.5
$ws.Range("L67").Value = 99999.5
$ws.Range("N67").Value = -101715.5
$ws.Range("H132").Value = 19447906
$ws.Range("I132").Value = 4184.8
$ws.Range("K132").Value = 12554.4
$ws.Range("M132").Value = -10024.4
